# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# Update the "OFF" sheet (row 3 = "R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 319
$wsOff.Range("C3").Value = 229
$wsOff.Range("D3").Value = 70
$wsOff.Range("E3").Value = 31
$wsOff.Range("G3").Value = 5

# Update the "DEF" sheet (row 3 = "R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 366
$wsDef.Range("C3").Value = 245
$wsDef.Range("D3").Value = 97
$wsDef.Range("E3").Value = 49
$wsDef.Range("F3").Value = 11
